$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New collapsed card rows (name + remaining fields as a Python-tuple-like string)
$ws.Range("A2").Value = "('Dromoka, the Eternal', ['{3}{G}{W}', 'Legendary Creature " + [char]0x2014 + " Dragon', 'Flying', 'Whenever a Dragon you control attacks, bolster 2. (Choose a creature with the least toughness among creatures you control and put two +1/+1 counters on it.)', '5/5'])"

$ws.Range("A3").Value = "('Honored Hierarch', ['{G}', 'Creature " + [char]0x2014 + " Human Druid', 'Renown 1 (When this creature deals combat damage to a player, if it isn" + [char]0x2019 + "t renowned, put a +1/+1 counter on it and it becomes renowned.)', 'As long as Honored Hierarch is renowned, it has vigilance and " + [char]0x201C + "{T}: Add one mana of any color." + [char]0x201D + "', '1/1'])"

$ws.Range("A4").Value = "('Sandsteppe Citadel', ['Land', 'Sandsteppe Citadel enters the battlefield tapped.', '{T}: Add {W}, {B}, or {G}.'])"

$ws.Range("A5").Value = "('Seeker of the Way', ['{1}{W}', 'Creature " + [char]0x2014 + " Human Warrior', 'Prowess (Whenever you cast a noncreature spell, this creature gets +1/+1 until end of turn.)', 'Whenever you cast a noncreature spell, Seeker of the Way gains lifelink until end of turn.', '2/2'])"

$ws.Range("A6").Value = "('Siege Rhino', ['{1}{W}{B}{G}', 'Creature " + [char]0x2014 + " Rhino', 'Trample', 'When Siege Rhino enters the battlefield, each opponent loses 3 life and you gain 3 life.', '4/5'])"

$ws.Range("A7").Value = "('Valorous Stance', ['{1}{W}', 'Instant', 'Choose one " + [char]0x2014 + "', '" + [char]0x2022 + " Target creature gains indestructible until end of turn.', '" + [char]0x2022 + " Destroy target creature with toughness 4 or greater.'])"

# Remove the now-unused trailing rows (8-35) so the sheet dimension shrinks to A1:A7
$ws.Rows("8:35").Delete()
